$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("YDS")

$ws.Range("B2").Value = "NCT(2.706298172812425, 1.2090313287454308, -0.0072667617676891155, 2.316990275100726)"
$ws.Range("C2").Value = "JSU(-1.3739278578640508, 1.18568690325675, 2.310451450803438, 4.622850275754908)"
$ws.Range("D2").Value = "NCT(2.4086928129537304, 1.1306078712125993, -0.016179999753329843, 2.4530510187780594)"
$ws.Range("E2").Value = "NIG(1.035249752960902, 0.8339133762786255, 3.9348229756215485, 4.340252631256927)"
